# Ajustar las categorías de densidad poblacional
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# States that were "High" density become "Very High"
$veryHighCells = @("D2", "D3")
foreach ($cell in $veryHighCells) {
    $ws.Range($cell).Value = "Very High"
}

# States that were "Low" density become "Very Low"
$veryLowCells = @("D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37")
foreach ($cell in $veryLowCells) {
    $ws.Range($cell).Value = "Very Low"
}

# Update the active selection to match the last edited cell
$ws.Range("E24").Select()
